# Mark attendance for 2025-05-20 (column U) across all four sheets.
# Sheet 1 (First_year): rows 2-3 => "A" (Absent), rows 4-18 => "P" (Present)
# Sheet 2 (Second_year): rows 2-10 => "P"
# Sheet 3 (Third_year): rows 2-23 => "P"
# Sheet 4 (Fourth_year): rows 2-31 => "P"

$wb = $excel.ActiveWorkbook

$sheetLastRows = @{
    "First_year"  = 18
    "Second_year" = 10
    "Third_year"  = 23
    "Fourth_year" = 31
}

$absentRows = @{
    "First_year" = @(2, 3)
}

foreach ($ws in $wb.Worksheets) {
    $name = $ws.Name
    if (-not $sheetLastRows.ContainsKey($name)) { continue }
    $lastRow = $sheetLastRows[$name]

    $absentSet = @()
    if ($absentRows.ContainsKey($name)) { $absentSet = $absentRows[$name] }

    for ($row = 2; $row -le $lastRow; $row++) {
        if ($absentSet -contains $row) {
            $ws.Range("U$row").Value = "A"
        } else {
            $ws.Range("U$row").Value = "P"
        }
    }
}
